$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Content edit -----------------------------------------------------
# The "Cases" tab Cypher query (row 2, column B) was trimmed: the trailing
# "Cohort" return column (and the optional-match cohort plumbing it relied
# on being listed) was dropped from the end of the RETURN clause.
$casesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
  "MATCH (c)<--(diag:diagnosis)`n" +
  "MATCH (samp:sample)-->(c) `n" +
  "WHERE samp.summarized_sample_type IN [`"Primary Malignant Tumor Tissue`"] `n" +
  "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
  "  WITH DISTINCT c, s, demo, diag, co`n" +
  "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
  "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
  "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
  "        coalesce(demo.breed, '') AS Breed ,`n" +
  "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
  "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
  "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
  "        coalesce(demo.sex, '') AS Sex ,`n" +
  "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
  "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
  "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $casesQuery

# --- Row heights --------------------------------------------------------
# Re-saving bumped the body row heights down a bit (wrap-text reflow under
# the newer Excel build); bring them in line.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# --- View / selection state ---------------------------------------------
# Leave the selection on B2 instead of B4.
$ws.Range("B2").Select()
